{"js": "// The underlying OOXML diff for this change is purely a re-serialisation of\n// existing XML attributes (namespace declarations on <w:document>, and\n// attribute order on elements such as <w:color>, <w:pgSz>, <w:pgMar>,\n// <w:rFonts>, <w:lang>, <w:latentStyles>/<w:lsdException>, <w:style>,\n// <w:tblInd> and <w:tblCellMar> inside word/document.xml and\n// word/styles.xml). Every \"-\"/\"+\" pair in the diff contains exactly the same\n// attribute names/values, just sorted alphabetically by the tool that\n// produced the diff - there is no textual, formatting, or structural change\n// to the document content itself (same paragraphs, same run text, same\n// field instructions, same colors/sizes/margins, same styles).\n//\n// Since attribute order inside an OOXML element is not semantically\n// meaningful (and is not something the Word JavaScript API exposes control\n// over \u2014 Office.js never lets you dictate XML attribute ordering), applying\n// this change means leaving the document's content/formatting untouched.\n// We simply touch the body to confirm the context/binding is valid, without\n// mutating anything.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying OOXML diff for this change is purely a re-serialisation of\n# existing XML attributes (namespace declarations on <w:document>, and\n# attribute order on elements such as <w:color>, <w:pgSz>, <w:pgMar>,\n# <w:rFonts>, <w:lang>, <w:latentStyles>/<w:lsdException>, <w:style>,\n# <w:tblInd> and <w:tblCellMar> inside word/document.xml and\n# word/styles.xml). Every \"-\"/\"+\" pair in the diff contains exactly the same\n# attribute names/values, just sorted alphabetically by the tool that\n# produced the diff - there is no textual, formatting, or structural change\n# to the document content itself (same paragraphs, same run text, same\n# field instructions, same colors/sizes/margins, same styles).\n#\n# Since attribute order inside an OOXML element is not semantically\n# meaningful (and is not something the Word object model exposes control\n# over \u2014 COM automation never lets you dictate XML attribute ordering),\n# applying this change means leaving the document's content/formatting\n# untouched. We simply touch the document to confirm the binding is valid,\n# without mutating anything.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
